# Updated cryptos list values (price + 1h volume change) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''68.649.98'
$ws.Range("E2").Value = '  +3.47%  '

$ws.Range("D3").Value = '''3.144.31'
$ws.Range("E3").Value = '  +2.47%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '''578.12'
$ws.Range("E5").Value = '  +0.64%  '

$ws.Range("D6").Value = '''179.59'
$ws.Range("E6").Value = '  +6.11%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").Value = '''3.143.00'
$ws.Range("E8").Value = '  +2.54%  '

$ws.Range("D9").Value = '''0.523'
$ws.Range("E9").Value = '  +2.81%  '

$ws.Range("E10").Value = '  +1.62%  '

$ws.Range("E11").Value = '  +2.15%  '

$ws.Range("E12").Value = '  +0.33%  '

$ws.Range("E13").Value = '  +1.83%  '

$ws.Range("D14").Value = '''36.96'
$ws.Range("E14").Value = '  +3.79%  '

$ws.Range("D15").Value = '''68.606.04'
$ws.Range("E15").Value = '  +3.38%  '

$ws.Range("D16").Value = '''3.670.14'
$ws.Range("E16").Value = '  +2.49%  '

$ws.Range("E17").Value = '  +0.90%  '

$ws.Range("D18").Value = '''7.16'
$ws.Range("E18").Value = '  +3.07%  '

$ws.Range("D19").Value = '''3.148.19'
$ws.Range("E19").Value = '  +2.50%  '

$ws.Range("E20").Value = '  -2.34%  '

$ws.Range("D21").Value = '''488.03'
$ws.Range("E21").Value = '  +0.10%  '

$ws.Range("D22").Value = '''0.699'
$ws.Range("E22").Value = '  +1.93%  '

$ws.Range("D23").Value = '''7.80'
$ws.Range("E23").Value = '  +0.82%  '

$ws.Range("D24").Value = '''84.06'
$ws.Range("E24").Value = '  +1.74%  '

$ws.Range("D25").Value = '''2.35'
$ws.Range("E25").Value = '  +6.97%  '

$ws.Range("D26").Value = '''13.04'
$ws.Range("E26").Value = '  +2.91%  '

$ws.Range("D27").Value = '''10.62'
$ws.Range("E27").Value = '  +4.91%  '

$ws.Range("E28").Value = '  +0.00%  '

$ws.Range("D29").Value = '''8.13'
$ws.Range("E29").Value = '  +4.37%  '

$ws.Range("D30").Value = '''2.36'
$ws.Range("E30").Value = '  +4.32%  '

$ws.Range("D31").Value = '''2.64'
$ws.Range("E31").Value = '  +1.55%  '

$ws.Range("D32").Value = '''28.17'
$ws.Range("E32").Value = '  +2.11%  '

$ws.Range("D33").Value = '''0.112'
$ws.Range("E33").Value = '  +1.31%  '

$ws.Range("D34").Value = '''0.0₃0949'
$ws.Range("E34").Value = '  +4.49%  '

$ws.Range("E35").Value = '  -0.05%  '

$ws.Range("E36").Value = '  +3.53%  '

$ws.Range("D37").Value = '''48.16'
$ws.Range("E37").Value = '  +2.45%  '

$ws.Range("D38").Value = '''0.961'
$ws.Range("E38").Value = '  +1.62%  '

$ws.Range("D39").Value = '''0.324'
$ws.Range("E39").Value = '  +7.96%  '

$ws.Range("D40").Value = '''2.04'
$ws.Range("E40").Value = '  +4.25%  '

$ws.Range("D41").Value = '''0.126'
$ws.Range("E41").Value = '  +3.18%  '

$ws.Range("D42").Value = '''49.21'
$ws.Range("E42").Value = '  +0.25%  '

$ws.Range("D43").Value = '''8.38'
$ws.Range("E43").Value = '  +1.17%  '

$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").Value = '''404.44'
$ws.Range("E44").Value = '  +10.58%  '

$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '''2.73'
$ws.Range("E45").Value = '  +7.96%  '

$ws.Range("D46").Value = '''27.88'
$ws.Range("E46").Value = '  +14.38%  '

$ws.Range("D47").Value = '''2.813.97'
$ws.Range("E47").Value = '  +1.52%  '

$ws.Range("D48").Value = '''0.0348'
$ws.Range("E48").Value = '  +1.16%  '

$ws.Range("D49").Value = '''135.35'
$ws.Range("E49").Value = '  +0.89%  '

$ws.Range("D51").Value = '''2.38'
$ws.Range("E51").Value = '  +10.42%  '

